# Update the 'F' column (numeric counter) values on each sheet to match the
# regenerated data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 544
$ws.Range("F3").Value = 970
$ws.Range("F4").Value = 66
$ws.Range("F7").Value = 1194
$ws.Range("F8").Value = 953
$ws.Range("F9").Value = 38
$ws.Range("F11").Value = 1053
$ws.Range("F12").Value = 3273
$ws.Range("F13").Value = 577
$ws.Range("F15").Value = 1703
$ws.Range("F17").Value = 649
$ws.Range("F18").Value = 21
$ws.Range("F20").Value = 378
$ws.Range("F23").Value = 773
$ws.Range("F24").Value = 660
$ws.Range("F25").Value = 517
$ws.Range("F26").Value = 485
$ws.Range("F28").Value = 49
$ws.Range("F29").Value = 1031
$ws.Range("F30").Value = 1166
$ws.Range("F31").Value = 336
$ws.Range("F34").Value = 1423
$ws.Range("F35").Value = 471
$ws.Range("F38").Value = 4049

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 197
$ws.Range("F8").Value = 11
$ws.Range("F9").Value = 8
$ws.Range("F12").Value = 398
$ws.Range("F24").Value = 1
$ws.Range("F37").Value = 17
$ws.Range("F38").Value = 10

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 99
$ws.Range("F3").Value = 50
$ws.Range("F4").Value = 1290
$ws.Range("F5").Value = 1687
$ws.Range("F6").Value = 460
$ws.Range("F7").Value = 1041
$ws.Range("F8").Value = 65

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1290
$ws.Range("F3").Value = 1687
$ws.Range("F4").Value = 460
$ws.Range("F5").Value = 1041
$ws.Range("F6").Value = 544
$ws.Range("F7").Value = 970
$ws.Range("F8").Value = 66
$ws.Range("F9").Value = 1194
$ws.Range("F10").Value = 953
$ws.Range("F12").Value = 38
$ws.Range("F14").Value = 197
$ws.Range("F15").Value = 197
$ws.Range("F16").Value = 11
$ws.Range("F17").Value = 1053
$ws.Range("F18").Value = 8
$ws.Range("F19").Value = 3299
$ws.Range("F20").Value = 577
$ws.Range("F22").Value = 1703
$ws.Range("F24").Value = 649
$ws.Range("F26").Value = 378
$ws.Range("F31").Value = 773
$ws.Range("F32").Value = 660
$ws.Range("F33").Value = 517
$ws.Range("F34").Value = 485
$ws.Range("F36").Value = 49
$ws.Range("F40").Value = 1031
$ws.Range("F41").Value = 1166
$ws.Range("F42").Value = 336
$ws.Range("F46").Value = 1423
$ws.Range("F47").Value = 471
$ws.Range("F50").Value = 4049
$ws.Range("F51").Value = 10
